$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for B1 and C1 (split the bracket text differently)
$ws.Range("B1").Value = "Ижорский (сойкинский) [Nirvi] @ Nirvi"
$ws.Range("C1").Value = "Ижорский (сойкинский) [учебная] @ Nirvi | в учебной системе"

# Row 1 height grew (likely due to wrapped text needing more vertical space)
$ws.Rows("1").RowHeight = 45

# Update the active selection from F1 to D1
$ws.Range("D1").Select()
